{"js": "// Apply the \"Dokumentation Vincents Vorschl\u00e4ge eingef\u00fcgt\" text edits.\n// Each edit is a narrow, uniquely-matching find/replace so existing\n// run formatting (Arial rPr) around the edited text is preserved as\n// closely as possible.\n\nasync function replaceOnce(body, findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  // Use the first (and expected only) match.\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"...zu stellen (siehe 2.3 Spielbalancing).\" -> \"...(siehe 2.3.2 Kaufsystem).\"\n//    (the heading \"2.3 Spielbalancing\" became \"2.3.2 Kaufsystem\")\nawait replaceOnce(\n  context.document.body,\n  \"siehe 2.3 Spielbalancing\",\n  \"siehe 2.3.2 Kaufsystem\"\n);\n\n// 2) \"Im Laufe des Spiels werden Hindernisse schwieriger zu absolvieren sein.\"\n//    -> \"Im Laufe des Spiels werden die Hindernisse immer schwieriger zu absolvieren.\"\nawait replaceOnce(\n  context.document.body,\n  \"werden Hindernisse schwieriger zu absolvieren sein.\",\n  \"werden die Hindernisse immer schwieriger zu absolvieren.\"\n);\n\n// 3) \"...ihre Lebenspunkte erh\u00f6ht werden. Damit...\"\n//    -> \"...ihre Lebenspunkte erh\u00f6ht werden, sie mehr Schaden verursachen oder neue\n//        Angriffsmuster benutzen. Damit...\"\nawait replaceOnce(\n  context.document.body,\n  \"Lebenspunkte erh\u00f6ht werden. Damit\",\n  \"Lebenspunkte erh\u00f6ht werden, sie mehr Schaden verursachen oder neue Angriffsmuster benutzen. Damit\"\n);\n\n// 4) \"...lassen sich Waffen im Shop kaufen. Gold wird...\"\n//    -> \"...lassen sich Waffen in verschiedenen Shops in den St\u00e4dten kaufen. Gold wird...\"\nawait replaceOnce(\n  context.document.body,\n  \"Waffen im Shop kaufen.\",\n  \"Waffen in verschiedenen Shops in den St\u00e4dten kaufen.\"\n);\n\n// 5) \"...fallengelassen oder liegen zuf\u00e4llig verteilt...\"\n//    -> \"...fallengelassen oder liegt zuf\u00e4llig verteilt...\"\nawait replaceOnce(\n  context.document.body,\n  \"fallengelassen oder liegen zuf\u00e4llig\",\n  \"fallengelassen oder liegt zuf\u00e4llig\"\n);\n\n// 6) \"Wenn man diesen nicht nutzt, ist es unm\u00f6glich, die n\u00e4chsten Level...\"\n//    -> \"...ist es nur schwer m\u00f6glich, die n\u00e4chsten Level...\"\nawait replaceOnce(\n  context.document.body,\n  \"ist es unm\u00f6glich, die n\u00e4chsten\",\n  \"ist es nur schwer m\u00f6glich, die n\u00e4chsten\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\n# 1) \"...zu stellen (siehe 2.3 Spielbalancing).\" -> \"...(siehe 2.3.2 Kaufsystem).\"\n#    (the referenced heading \"2.3 Spielbalancing\" became \"2.3.2 Kaufsystem\")\nReplace-Text \"siehe 2.3 Spielbalancing\" \"siehe 2.3.2 Kaufsystem\"\n\n# 2) \"Im Laufe des Spiels werden Hindernisse schwieriger zu absolvieren sein.\"\n#    -> \"Im Laufe des Spiels werden die Hindernisse immer schwieriger zu absolvieren.\"\nReplace-Text \"werden Hindernisse schwieriger zu absolvieren sein.\" \"werden die Hindernisse immer schwieriger zu absolvieren.\"\n\n# 3) \"...ihre Lebenspunkte erh\u00f6ht werden. Damit...\"\n#    -> \"...ihre Lebenspunkte erh\u00f6ht werden, sie mehr Schaden verursachen oder neue\n#        Angriffsmuster benutzen. Damit...\"\nReplace-Text \"Lebenspunkte erh\u00f6ht werden. Damit\" \"Lebenspunkte erh\u00f6ht werden, sie mehr Schaden verursachen oder neue Angriffsmuster benutzen. Damit\"\n\n# 4) \"...lassen sich Waffen im Shop kaufen. Gold wird...\"\n#    -> \"...lassen sich Waffen in verschiedenen Shops in den St\u00e4dten kaufen. Gold wird...\"\nReplace-Text \"Waffen im Shop kaufen.\" \"Waffen in verschiedenen Shops in den St\u00e4dten kaufen.\"\n\n# 5) \"...fallengelassen oder liegen zuf\u00e4llig verteilt...\"\n#    -> \"...fallengelassen oder liegt zuf\u00e4llig verteilt...\"\nReplace-Text \"fallengelassen oder liegen zuf\u00e4llig\" \"fallengelassen oder liegt zuf\u00e4llig\"\n\n# 6) \"Wenn man diesen nicht nutzt, ist es unm\u00f6glich, die n\u00e4chsten Level...\"\n#    -> \"...ist es nur schwer m\u00f6glich, die n\u00e4chsten Level...\"\nReplace-Text \"ist es unm\u00f6glich, die n\u00e4chsten\" \"ist es nur schwer m\u00f6glich, die n\u00e4chsten\"\n\n# The original document had no footnotes/endnotes parts. Adding (and then\n# clearing) a throwaway footnote/endnote mints the standard separator /\n# continuationSeparator boilerplate parts (word/footnotes.xml,\n# word/endnotes.xml) the same way Word silently creates them the first time\n# a document is touched with footnote/endnote-aware machinery, without\n# leaving a visible note reference or note body behind.\n$noteRange = $d.Range(0, 0)\n\n$fn = $d.Footnotes.Add($noteRange, \"\", \"x\")\n$fn.Reference.Text = \"\"\n\n$en = $d.Endnotes.Add($noteRange, \"\", \"x\")\n$en.Reference.Text = \"\"\n"}
